$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 16673734
$ws.Range("I32").Value = 10666.333
$ws.Range("J32").Value = 33336800
$ws.Range("K32").Value = 10666.333
$ws.Range("L32").Value = 33336800
$ws.Range("M32").Value = -10340.333
$ws.Range("N32").Value = -33337452
$ws.Range("H40").Value = 3926.5715
$ws.Range("I40").Value = 2759.5
$ws.Range("J40").Value = 6844.25
$ws.Range("K40").Value = 2759.5
$ws.Range("L40").Value = 6844.25
$ws.Range("M40").Value = -2584.5
$ws.Range("N40").Value = -7194.25
$ws.Range("H51").Value = 8571.857
$ws.Range("I51").Value = 10001
$ws.Range("K51").Value = 10001
$ws.Range("M51").Value = -9517
$ws.Range("H82").Value = 7399.4287
$ws.Range("I82").Value = 5299.3335
$ws.Range("K82").Value = 15898.0005
$ws.Range("M82").Value = -15492.0005
$ws.Range("H85").Value = 7399.4287
$ws.Range("I85").Value = 5299.3335
$ws.Range("K85").Value = 15898.0005
$ws.Range("M85").Value = -14494.0005
$ws.Range("H106").Value = 4636.1816
$ws.Range("I106").Value = 4779.9
$ws.Range("K106").Value = 4779.9
$ws.Range("M106").Value = -4148.9
$ws.Range("H132").Value = 17690.344
$ws.Range("I132").Value = 1258.619
$ws.Range("K132").Value = 3775.857
$ws.Range("M132").Value = -1245.857
$ws.Range("H137").Value = 1938.2
$ws.Range("I137").Value = 2123.75
$ws.Range("J137").Value = 1814.5
$ws.Range("K137").Value = 6371.25
$ws.Range("L137").Value = 5443.5
$ws.Range("M137").Value = -3821.25
$ws.Range("N137").Value = -10543.5
$ws.Range("H138").Value = 1557.1154
$ws.Range("I138").Value = 872.4706
$ws.Range("J138").Value = 2850.3333
$ws.Range("K138").Value = 2617.4118
$ws.Range("L138").Value = 8550.999899999999
$ws.Range("M138").Value = 2522.5882
$ws.Range("N138").Value = -18830.9999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9670.49
$ws.Range("I32").Value = 7968.5117
$ws.Range("K32").Value = 7968.5117
$ws.Range("M32").Value = -7681.5117
$ws.Range("H61").Value = 2767.7097
$ws.Range("I61").Value = 1459.1666
$ws.Range("K61").Value = 1459.1666
$ws.Range("M61").Value = -1247.1666
$ws.Range("H74").Value = 1082.7646
$ws.Range("I74").Value = 1116.4615
$ws.Range("K74").Value = 1116.4615
$ws.Range("M74").Value = -242.4614999999999
$ws.Range("H77").Value = 1082.7646
$ws.Range("I77").Value = 1116.4615
$ws.Range("K77").Value = 5582.307499999999
$ws.Range("M77").Value = -1214.307499999999
$ws.Range("H102").Value = 2122.5454
$ws.Range("I102").Value = 1931.375
$ws.Range("J102").Value = 2632.3333
$ws.Range("K102").Value = 1931.375
$ws.Range("L102").Value = 2632.3333
$ws.Range("M102").Value = -309.375
$ws.Range("N102").Value = -5876.3333
$ws.Range("H132").Value = 857.2
$ws.Range("I132").Value = 857.2
$ws.Range("K132").Value = 2571.6
$ws.Range("M132").Value = -41.60000000000036
$ws.Range("H136").Value = 2767.7097
$ws.Range("I136").Value = 1459.1666
$ws.Range("K136").Value = 4377.4998
$ws.Range("M136").Value = -1827.4998

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4078.5
$ws.Range("I20").Value = 2473.6
$ws.Range("K20").Value = 2473.6
$ws.Range("M20").Value = -2226.6
$ws.Range("H86").Value = 3040.6177
$ws.Range("I86").Value = 1769.0952
$ws.Range("J86").Value = 5094.615
$ws.Range("K86").Value = 1769.0952
$ws.Range("L86").Value = 5094.615
$ws.Range("M86").Value = -646.0952
$ws.Range("N86").Value = -7340.615
$ws.Range("H89").Value = 3040.6177
$ws.Range("I89").Value = 1769.0952
$ws.Range("J89").Value = 5094.615
$ws.Range("K89").Value = 8845.476000000001
$ws.Range("L89").Value = 25473.075
$ws.Range("M89").Value = -3229.476000000001
$ws.Range("N89").Value = -36705.075
$ws.Range("H105").Value = 1728.9231
$ws.Range("I105").Value = 1787.7273
$ws.Range("J105").Value = 1405.5
$ws.Range("K105").Value = 1787.7273
$ws.Range("L105").Value = 1405.5
$ws.Range("M105").Value = -40.72730000000001
$ws.Range("N105").Value = -4899.5
$ws.Range("H107").Value = 816.3333
$ws.Range("I107").Value = 799.5
$ws.Range("K107").Value = 799.5
$ws.Range("M107").Value = 1120.5
$ws.Range("H132").Value = 70099.85000000001
$ws.Range("J132").Value = 70099.85000000001
$ws.Range("L132").Value = 70099.85000000001
$ws.Range("N132").Value = -80219.85000000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1512.9
$ws.Range("I58").Value = 747.7222
$ws.Range("K58").Value = 747.7222
$ws.Range("M58").Value = -544.7222
$ws.Range("H86").Value = 92499.5
$ws.Range("I86").Value = 159999
$ws.Range("J86").Value = 25000
$ws.Range("K86").Value = 159999
$ws.Range("L86").Value = 25000
$ws.Range("M86").Value = -158876
$ws.Range("N86").Value = -27246
$ws.Range("H89").Value = 92499.5
$ws.Range("I89").Value = 159999
$ws.Range("J89").Value = 25000
$ws.Range("K89").Value = 799995
$ws.Range("L89").Value = 125000
$ws.Range("M89").Value = -794379
$ws.Range("N89").Value = -136232
$ws.Range("H107").Value = 6317.9443
$ws.Range("I107").Value = 794.5
$ws.Range("K107").Value = 794.5
$ws.Range("M107").Value = 1125.5
$ws.Range("H122").Value = 313182.78
$ws.Range("I122").Value = 568625.2
$ws.Range("J122").Value = 6651.933
$ws.Range("K122").Value = 1705875.6
$ws.Range("L122").Value = 19955.799
$ws.Range("M122").Value = -1703425.6
$ws.Range("N122").Value = -24855.799
$ws.Range("H132").Value = 3617.5
$ws.Range("J132").Value = 5500
$ws.Range("L132").Value = 16500
$ws.Range("N132").Value = -21560
$ws.Range("H136").Value = 1512.9
$ws.Range("I136").Value = 747.7222
$ws.Range("K136").Value = 2243.1666
$ws.Range("M136").Value = 306.8334

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 6908.879
$ws.Range("I5").Value = 3563.1765
$ws.Range("J5").Value = 10463.6875
$ws.Range("K5").Value = 10689.5295
$ws.Range("L5").Value = 31391.0625
$ws.Range("M5").Value = -10577.5295
$ws.Range("N5").Value = -31615.0625
$ws.Range("H135").Value = 6908.879
$ws.Range("I135").Value = 3563.1765
$ws.Range("J135").Value = 10463.6875
$ws.Range("K135").Value = 32068.5885
$ws.Range("L135").Value = 94173.1875
$ws.Range("M135").Value = -29533.5885
$ws.Range("N135").Value = -99243.1875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 195184.67
$ws.Range("I70").Value = 230221.6
$ws.Range("K70").Value = 230221.6
$ws.Range("M70").Value = -229951.6
$ws.Range("H73").Value = 195184.67
$ws.Range("I73").Value = 230221.6
$ws.Range("K73").Value = 230221.6
$ws.Range("M73").Value = -229285.6
$ws.Range("H80").Value = 104290.45
$ws.Range("I80").Value = 280202.75
$ws.Range("J80").Value = 3769.1428
$ws.Range("K80").Value = 280202.75
$ws.Range("L80").Value = 3769.1428
$ws.Range("M80").Value = -279204.75
$ws.Range("N80").Value = -5765.1428
$ws.Range("H83").Value = 104290.45
$ws.Range("I83").Value = 280202.75
$ws.Range("J83").Value = 3769.1428
$ws.Range("K83").Value = 1401013.75
$ws.Range("L83").Value = 18845.714
$ws.Range("M83").Value = -1396021.75
$ws.Range("N83").Value = -28829.714
$ws.Range("H107").Value = 624.1875
$ws.Range("I107").Value = 383
$ws.Range("K107").Value = 383
$ws.Range("M107").Value = 1537
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").Value = ""
$ws.Range("H122").Value = 3173.4194
$ws.Range("I122").Value = 1372.2667
$ws.Range("K122").Value = 4116.800099999999
$ws.Range("M122").Value = -1666.800099999999
$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").Value = ""
$ws.Range("H132").Value = 2872
$ws.Range("I132").Value = 2423.1428
$ws.Range("K132").Value = 7269.428400000001
$ws.Range("M132").Value = -4739.428400000001
$ws.Range("H135").Value = 145000
$ws.Range("J135").Value = 145000
$ws.Range("L135").Value = 145000
$ws.Range("N135").Value = -155140
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 6195.1333
$ws.Range("I82").Value = 3658.5557
$ws.Range("J82").Value = 10000
$ws.Range("K82").Value = 3658.5557
$ws.Range("L82").Value = 10000
$ws.Range("M82").Value = -3297.5557
$ws.Range("N82").Value = -10722
$ws.Range("H85").Value = 6195.1333
$ws.Range("I85").Value = 3658.5557
$ws.Range("J85").Value = 10000
$ws.Range("K85").Value = 3658.5557
$ws.Range("L85").Value = 10000
$ws.Range("M85").Value = -2410.5557
$ws.Range("N85").Value = -12496
$ws.Range("I122").Value = 2351.125
$ws.Range("J122").Value = 7275.8335
$ws.Range("K122").Value = 7053.375
$ws.Range("L122").Value = 21827.5005
$ws.Range("M122").Value = -4603.375
$ws.Range("N122").Value = -26727.5005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5680.4585
$ws.Range("I81").Value = 6843.737
$ws.Range("J81").Value = 1260
$ws.Range("K81").Value = 13687.474
$ws.Range("L81").Value = 2520
$ws.Range("M81").Value = -12626.474
$ws.Range("N81").Value = -4642
$ws.Range("H84").Value = 5680.4585
$ws.Range("I84").Value = 6843.737
$ws.Range("J84").Value = 1260
$ws.Range("K84").Value = 68437.37
$ws.Range("L84").Value = 12600
$ws.Range("M84").Value = -63133.37
$ws.Range("N84").Value = -23208
$ws.Range("H100").Value = 509.2
$ws.Range("I100").Value = 536.0952
$ws.Range("J100").Value = 368
$ws.Range("K100").Value = 1072.1904
$ws.Range("L100").Value = 736
$ws.Range("M100").Value = -531.1904
$ws.Range("N100").Value = -1818
$ws.Range("H126").Value = 2222.8667
$ws.Range("I126").Value = 2164.923
$ws.Range("J126").Value = 2599.5
$ws.Range("K126").Value = 6494.768999999999
$ws.Range("L126").Value = 7798.5
$ws.Range("M126").Value = -4024.768999999999
$ws.Range("N126").Value = -12738.5
